$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the row labels: A2 "LogisticAT" -> "LogisticIT", A3 "LogisticIT" -> "LogisticAT"
$ws.Range("A2").Value = "LogisticIT"
$ws.Range("A3").Value = "LogisticAT"

# Update row 2 values (new run with new features)
$ws.Range("B2").Value = 0.4409
$ws.Range("C2").Value = 0.4409
$ws.Range("D2").Value = 0.9217
$ws.Range("E2").Value = 0.7554
$ws.Range("F2").Value = 0.6881
$ws.Range("G2").Value = 0.6095
$ws.Range("H2").Value = 0.6094000000000001

# Update row 3 values
$ws.Range("B3").Value = 0.4533
$ws.Range("C3").Value = 0.4533
$ws.Range("D3").Value = 0.9162
$ws.Range("E3").Value = 0.7155
$ws.Range("F3").Value = 0.6629
$ws.Range("G3").Value = 0.6303
$ws.Range("H3").Value = 0.63
